$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.615.17'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.841.58'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('E4').Value = '  -2.24%  '
$ws.Range('D5').Value = '319.56'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('E6').Value = '  -1.94%  '
$ws.Range('E7').Value = '  -2.80%  '
$ws.Range('D8').Value = '0.3734'
$ws.Range('E8').Value = '  -1.83%  '
$ws.Range('D9').Value = '0.07328'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('D10').Value = '0.8731'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').Value = '21.52'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('D12').Value = '1.863.63'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '6.699'
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').Value = '5.433'
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('D15').Value = '0.07130'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').Value = '87.76'
$ws.Range('E16').Value = '  +3.92%  '
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').Value = '0.000008955'
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('E20').Value = '  -1.02%  '
$ws.Range('D21').Value = '27.632.58'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('D22').Value = '5.209'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('D23').Value = '11.08'
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('D24').Value = '2.085.05'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').Value = '2.010'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').Value = '155.82'
$ws.Range('E26').Value = '  -1.62%  '
$ws.Range('D27').Value = '18.54'
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('D28').Value = '2.140'
$ws.Range('E28').Value = '  +7.56%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '118.57'
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('D31').Value = '0.08935'
$ws.Range('E31').Value = '  -1.51%  '
$ws.Range('D32').Value = '1.220'
$ws.Range('E32').Value = '  -0.67%  '
$ws.Range('D33').Value = '0.7731'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').Value = '4.537'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('D35').Value = '2.887'
$ws.Range('E35').Value = '  -4.87%  '
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('D37').Value = '1.132'
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01971'
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05318'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').Value = '7.251'
$ws.Range('E40').Value = '  +4.98%  '
$ws.Range('D41').Value = '2.897'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -2.17%  '
$ws.Range('D43').Value = '0.1678'
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('D44').Value = '8.763'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '109.01'
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '10.64'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('D47').Value = '0.4726'
$ws.Range('D48').Value = '0.06484'
$ws.Range('E48').Value = '  -3.70%  '
$ws.Range('D49').Value = '1.015'
$ws.Range('E49').Value = '  -2.10%  '
$ws.Range('E50').Value = '  -1.77%  '
$ws.Range('E51').Value = '  -4.29%  '
